$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two trailing
#    spaces appended, followed by a new red-colored parenthetical comment
#    typed in three separate chunks (so it lands as three adjacent runs,
#    matching how the text was actually composed).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Collapse(0)
$r.InsertAfter("  ")
$r.Collapse(0)

$rA = $d.Range($r.End, $r.End)
$rA.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$rA.Font.Color = 255
$rA.Collapse(0)

$rB = $d.Range($rA.End, $rA.End)
$rB.InsertAfter("rsion for main branch")
$rB.Font.Color = 255
$rB.Collapse(0)

$rC = $d.Range($rB.End, $rB.End)
$rC.InsertAfter(")")
$rC.Font.Color = 255

# ---------------------------------------------------------------------------
# 2) Remove the trailing "ank God almighty, we are free at last." paragraph
#    (the final paragraph of the document, right after the Raven poem ends).
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$lastPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Drop the now-unused styles that were only there to support the removed
#    paragraph (and a couple of other never-applied leftovers). Deleting by
#    descending index keeps earlier indices valid while we work.
# ---------------------------------------------------------------------------
$styleNames = @(
    "podcast-tools__subscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading 4 Char",
    "Heading 2 Char",
    "Hyperlink",
    "apple-converted-space",
    "Heading 4",
    "Heading 2"
)

$indices = @()
foreach ($nm in $styleNames) {
    $s = $d.Styles.Item($nm)
    $indices += $s.Index
}
$indices = $indices | Sort-Object -Descending

foreach ($idx in $indices) {
    $d.Styles.Item($idx).Delete()
}

Write-Host "Edit complete"
